$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 307.33334
$ws.Range("J9").Value = 307.33334
$ws.Range("L9").Value = 307.33334
$ws.Range("N9").Value = -645.33334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2498.4
$ws.Range("I32").Value = 3000
$ws.Range("J32").Value = 1746
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 1746
$ws.Range("M32").Value = -2674
$ws.Range("N32").Value = -2398

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1463.6666
$ws.Range("I62").Value = 1445
$ws.Range("K62").Value = 1445
$ws.Range("M62").Value = -821

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3833.5557
$ws.Range("I64").Value = 3833.5557
$ws.Range("K64").Value = 3833.5557
$ws.Range("M64").Value = -3585.5557

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 1463.6666
$ws.Range("I65").Value = 1445
$ws.Range("K65").Value = 7225
$ws.Range("M65").Value = -4105

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3833.5557
$ws.Range("I67").Value = 3833.5557
$ws.Range("K67").Value = 3833.5557
$ws.Range("M67").Value = -2975.5557

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 516.6667
$ws.Range("I129").Value = 516.6667
$ws.Range("K129").Value = 1550.0001
$ws.Range("M129").Value = 3449.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 531.5714
$ws.Range("I5").Value = 360.25
$ws.Range("J5").Value = 760
$ws.Range("K5").Value = 360.25
$ws.Range("L5").Value = 760
$ws.Range("M5").Value = -248.25
$ws.Range("N5").Value = -984

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3963.625
$ws.Range("I32").Value = 3963.625
$ws.Range("K32").Value = 3963.625
$ws.Range("M32").Value = -3676.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3369.9
$ws.Range("I45").Value = 1199.5
$ws.Range("J45").Value = 3912.5
$ws.Range("K45").Value = 1199.5
$ws.Range("L45").Value = 3912.5
$ws.Range("M45").Value = -822.5
$ws.Range("N45").Value = -4666.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1895.75
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1895.75
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1093.5454
$ws.Range("I97").Value = 985.5333000000001
$ws.Range("J97").Value = 1325
$ws.Range("K97").Value = 985.5333000000001
$ws.Range("L97").Value = 1325
$ws.Range("M97").Value = -489.5333000000001
$ws.Range("N97").Value = -2317

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 531.5714
$ws.Range("I4").Value = 360.25
$ws.Range("J4").Value = 760
$ws.Range("K4").Value = 360.25
$ws.Range("L4").Value = 760
$ws.Range("M4").Value = -245.25
$ws.Range("N4").Value = -990

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2091.375
$ws.Range("J94").Value = 1215.5714
$ws.Range("L94").Value = 1215.5714
$ws.Range("N94").Value = -2117.5714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 2666.3333
$ws.Range("I15").Value = 2999.5
$ws.Range("K15").Value = 2999.5
$ws.Range("M15").Value = -2829.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1718.5
$ws.Range("I58").Value = 1966.3334
$ws.Range("K58").Value = 1966.3334
$ws.Range("M58").Value = -1763.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4193.6
$ws.Range("I132").Value = 3994.5
$ws.Range("K132").Value = 11983.5
$ws.Range("M132").Value = -9453.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1718.5
$ws.Range("I136").Value = 1966.3334
$ws.Range("K136").Value = 5899.0002
$ws.Range("M136").Value = -3349.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 31.9
$ws.Range("I2").Value = 18.428572
$ws.Range("J2").Value = 63.333332
$ws.Range("K2").Value = 110.571432
$ws.Range("L2").Value = 379.999992
$ws.Range("M2").Value = 2.428568000000013
$ws.Range("N2").Value = -605.999992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 150
$ws.Range("I31").Value = 150
$ws.Range("K31").Value = 450
$ws.Range("M31").Value = -162

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 769.4
$ws.Range("J34").Value = 849.25
$ws.Range("L34").Value = 2547.75
$ws.Range("N34").Value = -2715.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 910.5217
$ws.Range("I131").Value = 517
$ws.Range("J131").Value = 1049.4117
$ws.Range("K131").Value = 1551
$ws.Range("L131").Value = 3148.2351
$ws.Range("M131").Value = 3489
$ws.Range("N131").Value = -13228.2351

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3608.611
$ws.Range("I80").Value = 1789.6
$ws.Range("J80").Value = 4308.231
$ws.Range("K80").Value = 1789.6
$ws.Range("L80").Value = 4308.231
$ws.Range("M80").Value = -791.5999999999999
$ws.Range("N80").Value = -6304.231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3608.611
$ws.Range("I83").Value = 1789.6
$ws.Range("J83").Value = 4308.231
$ws.Range("K83").Value = 8948
$ws.Range("L83").Value = 21541.155
$ws.Range("M83").Value = -3956
$ws.Range("N83").Value = -31525.155

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H120").Value = 45000
$ws.Range("I120").Value = 45000
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 45000
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -40162
$ws.Range("N120").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 7666.6665
$ws.Range("I126").Value = 4000
$ws.Range("K126").Value = 12000
$ws.Range("M126").Value = -9530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 3500
$ws.Range("I45").Value = 3500
$ws.Range("K45").Value = 3500
$ws.Range("M45").Value = -3093

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 16000
$ws.Range("I48").Value = 16000
$ws.Range("J48").Value = 16000
$ws.Range("K48").Value = 16000
$ws.Range("L48").Value = 16000
$ws.Range("M48").Value = -15339
$ws.Range("N48").Value = -17322

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 9000
$ws.Range("J98").Value = 9000
$ws.Range("L98").Value = 9000
$ws.Range("N98").Value = -14990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1342.5714
$ws.Range("I96").Value = 1149.75
$ws.Range("J96").Value = 1599.6666
$ws.Range("K96").Value = 1149.75
$ws.Range("L96").Value = 1599.6666
$ws.Range("M96").Value = 223.25
$ws.Range("N96").Value = -4345.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1195.6
$ws.Range("J107").Value = 1799.5
$ws.Range("L107").Value = 5398.5
$ws.Range("N107").Value = -9238.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3197.4666
$ws.Range("I122").Value = 3397.4546
$ws.Range("K122").Value = 10192.3638
$ws.Range("M122").Value = -7742.363799999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1154.8
$ws.Range("I126").Value = 1154.8
$ws.Range("K126").Value = 3464.4
$ws.Range("M126").Value = -994.3999999999996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1943.7693
$ws.Range("I132").Value = 1943.7693
$ws.Range("M132").Value = -3301.3079
